$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set new column widths for F and G (closest achievable values; engine
# quantizes stored width to pixel granularity, so inputs are tuned to land
# nearest the target 20.85546875 / 27.42578125 column widths) ---
$ws.Columns.Item(6).ColumnWidth = 20
$ws.Columns.Item(7).ColumnWidth = 26.6

# --- Set row heights ---
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 32.25

# --- Populate the new WNBA sponsorship comparison table (rows 1-11) and shared columns F/G ---
# Row 1
$ws.Range('A1').Value = 'Tenor of WNBA Sponserships'
$ws.Range('B1').Value = 'Main Sponser''s Name'
$ws.Range('C1').Value = 'WNBA Team being Sponsored'
$ws.Range('D1').Value = 'estimated Sponsership Amount'
$ws.Range('E1').Value = 'per year'
$ws.Range('F1').Value = 'Shared Sponserships'
$ws.Range('G1').Value = 'Team Sponser is Shared With'

# Row 2
$ws.Range('A2').Value = '2021 N/A'
$ws.Range('B2').Value = 'Gateway Center Arena'
$ws.Range('C2').Value = 'Atlanta Dreams'
$ws.Range('D2').Value = 'N/A'
$ws.Range('E2').Value = 'N/A'
$ws.Range('F2').Value = 'No'
$ws.Range('G2').Value = 'not shared'

# Row 3
$ws.Range('A3').Value = '2016 15 year deal'
$ws.Range('B3').Value = 'Wintrust Financial'
$ws.Range('C3').Value = 'Chicago Sky'
$ws.Range('D3').Value = 'N/A'
$ws.Range('E3').Value = 'N/A'
$ws.Range('F3').Value = 'No'
$ws.Range('G3').Value = 'not shared'

# Row 4
$ws.Range('A4').Value = '2015 N/A'
$ws.Range('B4').Value = 'College Park Center'
$ws.Range('C4').Value = 'Dallas Wings'
$ws.Range('D4').Value = 'N/A'
$ws.Range('E4').Value = 'N/A'
$ws.Range('F4').Value = 'No'
$ws.Range('G4').Value = 'not shared'

# Row 5
$ws.Range('A5').Value = '2021 multi year deal'
$ws.Range('B5').Value = 'Gainbridge'
$ws.Range('C5').Value = 'Indiana Fever'
$ws.Range('D5').Value = 'N/A'
$ws.Range('E5').Value = 'N/A'
$ws.Range('F5').Value = 'Yes'
$ws.Range('G5').Value = 'Indiana Pacers'

# Row 6
$ws.Range('A6').Value = '2001 start 20 years'
$ws.Range('B6').Value = 'crypto.com center'
$ws.Range('C6').Value = 'Los Angeles Sparks'
$ws.Range('D6').Value = '700 mil total'
$ws.Range('E6').Value = '35 mil per year'
$ws.Range('F6').Value = 'Yes'
$ws.Range('G6').Value = 'Los Angeles Lakers'

# Row 7
$ws.Range('A7').Value = '2021 multi year deal'
$ws.Range('B7').Value = 'Michelob Ultra'
$ws.Range('C7').Value = 'Las Vegas Aces'
$ws.Range('D7').Value = 'N/A'
$ws.Range('E7').Value = 'N/A'
$ws.Range('F7').Value = 'No'
$ws.Range('G7').Value = 'not shared'

# Row 8
$ws.Range('A8').Value = '1990 - present (3 to 5 year
 annul renewel)'
$ws.Range('B8').Value = 'Target'
$ws.Range('C8').Value = 'Minesota Lynx'
$ws.Range('D8').Value = 'N/A'
$ws.Range('E8').Value = '1.5 mil per year'
$ws.Range('F8').Value = 'Yes'
$ws.Range('G8').Value = 'Minnesota Timberwolves'

# Row 9
$ws.Range('A9').Value = '2020 N/A'
$ws.Range('B9').Value = 'Barclays Center'
$ws.Range('C9').Value = 'New York Liberty'
$ws.Range('D9').Value = 'N/A'
$ws.Range('E9').Value = 'N/A'
$ws.Range('F9').Value = 'No'
$ws.Range('G9').Value = 'not shared'

# Row 10
$ws.Range('A10').Value = '2021 start N/A'
$ws.Range('B10').Value = 'Footprint Center'
$ws.Range('C10').Value = 'Phoenix Mercury'
$ws.Range('D10').Value = 'N/A'
$ws.Range('E10').Value = 'N/A'
$ws.Range('F10').Value = 'Yes'
$ws.Range('G10').Value = 'Phoenix Suns'

# Row 11
$ws.Range('A11').Value = '2020 start N/A'
$ws.Range('B11').Value = 'Amazon/ Climate 
Pledge Arena'
$ws.Range('C11').Value = 'Seattle Storm'
$ws.Range('D11').Value = 'N/A'
$ws.Range('E11').Value = 'N/A'
$ws.Range('F11').Value = 'No'
$ws.Range('G11').Value = 'not shared'

# Row 13
$ws.Range('A13').Value = 'Tenor of NBA sponserships'
$ws.Range('B13').Value = 'Main Sponser''s Name'
$ws.Range('C13').Value = 'NBA Teams Sponsered'
$ws.Range('D13').Value = 'estimated Sponsership Amount'
$ws.Range('E13').Value = 'per year'

# Row 14
$ws.Range('A14').Value = '2018 start 20 years'
$ws.Range('B14').Value = 'State Farm'
$ws.Range('C14').Value = 'Atlanta Hawks'
$ws.Range('D14').Value = '175 mil'
$ws.Range('E14').Value = '8.75 mil per year'

# Row 15
$ws.Range('A15').Value = '1994 - present (renewed in 
2013 for additional 20 years)'
$ws.Range('B15').Value = 'United Airlines
/United Center Joint Venture'
$ws.Range('C15').Value = 'Chicago Bulls'
$ws.Range('D15').Value = 'N/A'
$ws.Range('E15').Value = 'N/A'

# Row 16
$ws.Range('A16').Value = '1999 start 30 years'
$ws.Range('B16').Value = 'American Airlines'
$ws.Range('C16').Value = 'Dallas Mavricks'
$ws.Range('D16').Value = '195 mil'
$ws.Range('E16').Value = '6.5 mil per year'

# Row 17
$ws.Range('A17').Value = '2019  start 20 years'
$ws.Range('B17').Value = 'chase'
$ws.Range('C17').Value = 'Golden State Warriors'
$ws.Range('D17').Value = '300 mil total'
$ws.Range('E17').Value = '15 mil per year'

# Row 18
$ws.Range('A18').Value = '2021 multi year deal'
$ws.Range('B18').Value = 'Gainbridge Insurance Agency LLC /Group One
 Thousand One LLC'
$ws.Range('C18').Value = 'Indiana Pacers'
$ws.Range('D18').Value = 'N/A'
$ws.Range('E18').Value = 'N/A'

# Row 19
$ws.Range('A19').Value = '2021 start 20 years'
$ws.Range('B19').Value = 'crypto.com'
$ws.Range('C19').Value = 'Los Angeles Lakers'
$ws.Range('D19').Value = '700 mil total'
$ws.Range('E19').Value = '35 mil per year'

# Row 20
$ws.Range('A20').Value = '1990 - present (3 to 5 year
 annul renewel)'
$ws.Range('B20').Value = 'Target'
$ws.Range('C20').Value = 'Minnesota Timberwolves'
$ws.Range('D20').Value = 'N/A'
$ws.Range('E20').Value = '1.5 mil per year'

# Row 21
$ws.Range('A21').Value = '2021 NA'
$ws.Range('B21').Value = 'Madison Square Garden Sports / New York City'
$ws.Range('C21').Value = 'New York Knicks'
$ws.Range('D21').Value = 'N/A'
$ws.Range('E21').Value = 'N/A'

# Row 22
$ws.Range('A22').Value = '2021 start 15 years'
$ws.Range('B22').Value = 'Paycom'
$ws.Range('C22').Value = 'Oklahoma City Thunder'
$ws.Range('D22').Value = 'N/A'
$ws.Range('E22').Value = 'N/A'

# Row 23
$ws.Range('A23').Value = '2021 start N/A'
$ws.Range('B23').Value = 'Footprint'
$ws.Range('C23').Value = 'Phoenix Suns'
$ws.Range('D23').Value = 'N/A'
$ws.Range('E23').Value = 'N/A'

# --- Apply wrap-text formatting to cells that need it ---
$ws.Range('A8').WrapText = $true
$ws.Range('B11').WrapText = $true
$ws.Range('A15').WrapText = $true
$ws.Range('B15').WrapText = $true
$ws.Range('B18').WrapText = $true
$ws.Range('A20').WrapText = $true
$ws.Range('B21').WrapText = $true

# --- Update selection to match the final saved state ---
$ws.Range('F5').Select()
